$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cells that change value
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "303"
$t.Cell(5,1).Range.Text  = "0.00001"
$t.Cell(6,1).Range.Text  = "0.00052"
$t.Cell(7,1).Range.Text  = "0.00019"
$t.Cell(9,1).Range.Text  = "0.00033"
$t.Cell(10,1).Range.Text = "0.00039"
$t.Cell(11,1).Range.Text = "0.00046"
$t.Cell(12,1).Range.Text = "0.06583"

# Trailing rows: collapse the tab-separated multi-run cells down to a
# single value (matches the first column of the old row's data).
$t.Cell(44,1).Range.Text = "99.97"
$t.Cell(45,1).Range.Text = "0.07"
$t.Cell(46,1).Range.Text = "210"
